$wb = $excel.ActiveWorkbook

# Sheets that contain the data rows needing updates: "展览" and "全部类型"
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 523
    $ws.Range("F4").Value = 167
    $ws.Range("F6").Value = 87
    $ws.Range("F7").Value = 737
    $ws.Range("F8").Value = 4
    $ws.Range("F9").Value = 415
}
